# Update "Correspond Handoff Datetime" (D5) and "Correspond Handback DateTime" (G5)
# timestamps on the "zh-cn" and "de-de" worksheets, simulating a freshly
# generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-22 05:15:50"
$wsZhCn.Range("G5").Value = "2016-02-22 05:16:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-22 05:16:04"
$wsDeDe.Range("G5").Value = "2016-02-22 05:17:01"
